# The SmartCampus "room" import template used to be a verbose 6-column
# sheet (EN/KA header pairs explaining how to fill every numeric code).
# It is simplified down to a compact 3-column entry template:
#   A1 = type, B1 = room number, C1 = student capacity.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (the old second header / example row) is no longer needed ---
$ws.Range("A2:F2").Clear()

# --- Columns E and F of row 1 are no longer used either ---
$ws.Range("E1:F1").Clear()

# --- Row 1: new, shorter header labels ---
$ws.Range("A1").Value = "ტიპი"
$ws.Range("B1").Value = "ოთახის N"
$ws.Range("C1").Value = "სტუდ. ტევადობა"

# B1 carries rich-text formatting: "ოთახის" in bold Menlo, " N" regular Calibri
$ws.Range("B1").Characters(1, 6).Font.Name = "Menlo Bold"
$ws.Range("B1").Characters(1, 6).Font.Size = 12
$ws.Range("B1").Characters(7, 2).Font.Name = "Calibri"
$ws.Range("B1").Characters(7, 2).Font.Size = 11

# B1/C1/D1 drop the old numeric-format column style (back to Normal/General)
$ws.Range("B1:D1").Style = "Normal"
$ws.Range("D1").ClearContents()
$ws.Range("D1").Style = "Normal"

# Row 1 is now a bit taller to fit the two-line-ish header text
$ws.Rows(1).RowHeight = 15.75

# Move the active selection to the first data entry cell
$ws.Range("A2").Select()
